# Update weekly portfolio figures (Weekly Performance %, Performance %, Value £)
# for rows 2-18 on the active sheet, reflecting the latest recalculated data dump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.51, -36.61, 658.4213786315918)
    3  = @(0.09, 11.15, 1237.782521437694)
    4  = @(0.49, -39.98, 244.9175742082712)
    5  = @(-0.39, -22.16, 199.1211951205687)
    6  = @(2.04, 139.73, 1346.135633020761)
    7  = @(1.86, -26.28, 533.99789997018)
    8  = @(2.54, -31.29, 316.0115210622402)
    9  = @(8.77, -64.26000000000001, 431.4000091552734)
    10 = @(-0.7, -2.13, 591.9)
    11 = @(0.8, 1.96, 521.2191378822326)
    12 = @(-2.15, 45.26, 867.966730676651)
    13 = @(-1.28, -6.21, 1392.426816291997)
    14 = @(-0.5600000000000001, 30.6, 913.414410029291)
    15 = @(-1.97, -32.95, 465.3151275135048)
    16 = @(0.3, 26.28, 566.33290277462)
    17 = @(4.48, -72.03, 106.6252210095424)
    18 = @(0.33, 3.93, 10392.98807878442)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
}
